# issue #5: stock data from json to db
# Add three new columns (category, source_file, index) to the "股票" (stock)
# worksheet, shifting the former legislator_id column to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(5)

# --- copy existing cell formatting into the new cells -----------------
# Header row (row 1) uses the bold/centered/bordered style found on K1.
$ws.Range("K1").Copy()
$ws.Range("L1:N1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data row (row 2) uses the plain style found on K2.
$ws.Range("K2").Copy()
$ws.Range("L2:N2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- header row (row 1): insert "category" before "date" and append ---
# --- "source_file" / "index" at the end --------------------------------
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# --- data row (row 2): fill in the corresponding values ----------------
$ws.Range("I2").Value = "normal"

# Force J2 ("2012-03-26") to be stored as text rather than being
# auto-parsed into a date serial number, then restore the plain
# (non-bordered) data-row formatting used by the rest of row 2.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "2012-03-26"
$ws.Range("K2").Copy()
$ws.Range("J2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("K2").Value = "潘維剛"
$ws.Range("L2").Value = 678
$ws.Range("M2").Value = "tmp71a01"
$ws.Range("N2").Value = 71

Write-Output "stock sheet updated"
